# Case_3_229 (380 kV) bus voltage magnitude results (res_bus/vm_pu.xlsx)
# Updated slack/bus setpoint (column B) from 1.05 pu to 1.02 pu and refreshed the
# resulting per-bus voltage magnitudes (columns C:F, I:N) for rows 2-25 (buses 0-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (B:F, I:N -- column G stays 1.0 and column H is unused)
$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")
$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# One row of new vm_pu values per bus (row 2 = bus 0 ... row 25 = bus 23),
# in the same column order as $columns above.
$newValues = @(
    @(1.02, 1.023188248718674, 1.028624566296347, 1.023841431239618, 1.037601164883557, 1.033721375922591, 1.0283698706378, 1.031440840929838, 1.02667168423602, 1.040391583985057, 1.013520043920996),
    @(1.02, 1.023947941862817, 1.029197239628194, 1.024481170809434, 1.03976245779677, 1.033957088156699, 1.028768919443265, 1.031822123701145, 1.02711884403465, 1.042359146019254, 1.013652721290149),
    @(1.02, 1.0244398088786, 1.029567874401122, 1.024895773767951, 1.041155598438508, 1.034108016894516, 1.029026750356593, 1.032068214857898, 1.02740814811355, 1.043626594694832, 1.013738424989113),
    @(1.02, 1.024646659682289, 1.029723707240457, 1.025070227245104, 1.041740022020197, 1.034171087158562, 1.029135051274282, 1.032171522220882, 1.027529761988014, 1.044158095010291, 1.013774419448877),
    @(1.02, 1.024681394922618, 1.029749873301767, 1.025099527778654, 1.041838076729413, 1.03418165467256, 1.029153230122506, 1.032188859225603, 1.027550180919715, 1.044247258778881, 1.013780461000183),
    @(1.02, 1.024442572555322, 1.029569956577074, 1.024898104218367, 1.041163412420287, 1.034108861134931, 1.029028197837323, 1.032069595842747, 1.027409773162215, 1.043633701839595, 1.013738906087961),
    @(1.02, 1.023444929155903, 1.02881808778072, 1.024057499373395, 1.038332712351843, 1.033801366485383, 1.02850480972366, 1.031569826556983, 1.026822811653074, 1.041057727296836, 1.013564913339453),
    @(1.02, 1.021689239070095, 1.027493803962557, 1.022581252982988, 1.033302219466217, 1.033247269282378, 1.027579625016719, 1.030684382036707, 1.02578823491012, 1.03647362293883, 1.013257189859913),
    @(1.02, 1.020520349361309, 1.026611373712346, 1.021600504634179, 1.029918074003136, 1.032869556319939, 1.026960883917485, 1.030090855569524, 1.025098352331902, 1.033385544254663, 1.013051287051237),
    @(1.02, 1.020014585837691, 1.026229375468611, 1.021176649694822, 1.02844503252302, 1.03270401226716, 1.026692499258833, 1.029833082347137, 1.024799590021385, 1.032040371433604, 1.012961950450123),
    @(1.02, 1.019826779249529, 1.026087499443394, 1.021019334349074, 1.027896686768256, 1.03264222094304, 1.026592739182007, 1.029737217454959, 1.024688610807868, 1.031539474796041, 1.012928739917834),
    @(1.02, 1.019867061837178, 1.026117931646066, 1.021053073436108, 1.028014363378594, 1.03265548903459, 1.026614141229615, 1.029757786075726, 1.024712416475689, 1.031646975493923, 1.012935864908959),
    @(1.02, 1.019999060531659, 1.026217647642848, 1.02116364344092, 1.028399730642563, 1.032698910724778, 1.026684254489747, 1.029825160503039, 1.024790416557864, 1.031998992655409, 1.012959205806779),
    @(1.02, 1.020080396782592, 1.026279087967636, 1.021231785623176, 1.028637008955424, 1.032725624340299, 1.026727444295583, 1.029866656669275, 1.024838474248088, 1.032215716720927, 1.012973583334208),
    @(1.02, 1.020553923123214, 1.026636727793544, 1.021628651735826, 1.030015669780842, 1.032880500828605, 1.026978685903962, 1.03010794682316, 1.025118179407737, 1.033474647075624, 1.013057212249305),
    @(1.02, 1.020851053955156, 1.026861092647278, 1.021877814231243, 1.030878383329058, 1.032977116339543, 1.027136158387723, 1.030259094707699, 1.025293620885983, 1.034262170288529, 1.013109622450852),
    @(1.02, 1.021024401349095, 1.026991970501764, 1.022023225143729, 1.031380850741338, 1.033033278461462, 1.027227964500971, 1.03034718222109, 1.025395949195465, 1.034720748411305, 1.01314017513301),
    @(1.02, 1.021083514425116, 1.027036598096858, 1.022072819823926, 1.031552054973454, 1.033052395763902, 1.027259260395903, 1.030377205163369, 1.025430839863762, 1.034876982017002, 1.01315058987315),
    @(1.02, 1.020819170901781, 1.02683701942476, 1.021851073313829, 1.0307858990219, 1.032966770288195, 1.027119267727809, 1.030242885677964, 1.025274798060736, 1.034177756487909, 1.013104001121831),
    @(1.02, 1.019960188643265, 1.026188283348632, 1.021131079922406, 1.028286282818125, 1.032686132438746, 1.026663609817242, 1.029805323631174, 1.02476744763048, 1.031895366881042, 1.012952333237866),
    @(1.02, 1.019420439494478, 1.025780484807849, 1.020679104928604, 1.026707755755178, 1.032507942572226, 1.026376714160805, 1.02952953709588, 1.024448424027877, 1.030453148221941, 1.012856817741459),
    @(1.02, 1.01970653959961, 1.025996658195928, 1.020918637565975, 1.027545231565467, 1.032602570044467, 1.026528841448428, 1.029675800778075, 1.024617547468671, 1.031218389132226, 1.012907467090733),
    @(1.02, 1.020833577370408, 1.026847897048373, 1.021863156140748, 1.030827690980802, 1.032971445817027, 1.02712690002902, 1.030250210071179, 1.025283303297067, 1.034215901831398, 1.013106541211993),
    @(1.02, 1.02214285265814, 1.027836089312009, 1.022962299527358, 1.03460795045666, 1.033391976456218, 1.027819151566173, 1.030913859054774, 1.026055728257127, 1.037664230731661, 1.013336876830971)
)

for ($r = 0; $r -lt $newValues.Length; $r++) {
    $rowNum = 2 + $r
    $rowVals = $newValues[$r]
    for ($c = 0; $c -lt $colIndexes.Length; $c++) {
        $ws.Cells.Item($rowNum, $colIndexes[$c]).Value = $rowVals[$c]
    }
}
